$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B2: Programs tab query - now asks for "Special Topic" instead of "Focus Area",
# and both CASE branches of the Data Location Details column use program_acronym.
$programQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Special Topic",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.program_acronym     
        ELSE prg.data_link
    END AS "Data Location Details" 
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Lung Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

# B3: Projects tab query - org_name renamed to project_org_name.
$projectsQuery = @'
SELECT DISTINCT
    prj.project_id AS "Project ID", 
    prj.project_title AS "Project Title",
    prj.project_org_name AS "Organization",
    prj.project_start_date AS "Project Start Date",
    prj.project_end_date AS "Project End Date"
FROM 
    df_project prj
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Lung Cancer%'
ORDER BY 
    lower(prj.project_id) ASC
LIMIT 100;
'@

# B4: Grants tab query - project_end_date renamed to grant_end_date, extra space before LIKE.
$grantsQuery = @'
SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.grant_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type  LIKE '%Lung Cancer%'
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;
'@

# B5: Publications tab query - pub.title renamed to pub.publication_title, added a
# relative_citation_ratio = 1.0 case, extra space before LIKE.
$publicationsQuery = @'
SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.publication_title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN '0'
    WHEN pub.relative_citation_ratio = 7.0 THEN '7'
    WHEN pub.relative_citation_ratio = 2.0 THEN '2'
  WHEN pub.relative_citation_ratio = 1.0 THEN '1'
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
     prg.cancer_type  LIKE '%Lung Cancer%'
ORDER BY 
    lower(pub.pmid) ASC
LIMIT 100;
'@

# Update the query cells in place (existing cell formatting - wrap text, font - is kept).
$ws.Cells.Item(2, 2).Value = $programQuery
$ws.Cells.Item(3, 2).Value = $projectsQuery
$ws.Cells.Item(4, 2).Value = $grantsQuery
$ws.Cells.Item(5, 2).Value = $publicationsQuery

# Re-apply the shared "wrap text, size 12" formatting to the edited cells (and to C2,
# whose query text did not change but whose formatting was refreshed alongside the rest).
$ws.Cells.Item(2, 2).Font.Size = 12
$ws.Cells.Item(2, 2).WrapText = $true
$ws.Cells.Item(2, 3).Font.Size = 12
$ws.Cells.Item(2, 3).WrapText = $true
$ws.Cells.Item(3, 2).Font.Size = 12
$ws.Cells.Item(3, 2).WrapText = $true
$ws.Cells.Item(4, 2).Font.Size = 12
$ws.Cells.Item(4, 2).WrapText = $true
$ws.Cells.Item(5, 2).Font.Size = 12
$ws.Cells.Item(5, 2).WrapText = $true

# Move the active selection to B2.
$ws.Cells.Item(2, 2).Select()
